$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(1, "first-year", 0, "True: In Java, every variable has a type that is known at compile time.", "Correct", 1)
    ,@(2, "first-year", 0, "False: In Java, the result of applying one of the arithmetic operators to two double operands may produce a run-time exception if the operands cause an overflow or underflow.", "Incorrect", 1)
    ,@(3, "first-year", 0, "False: In Java, if you attempt to use a local variable of type int in an expression before that variable has been assigned a value, Java will throw a compile-time error.", "Correct", 1)
    ,@(4, "first-year", 0, "True: In Java, if a variable is declared and initialized in the body of a for loop, that variable cannot be accessed outside that loop.", "Correct", 1)
    ,@(5, "first-year", 0, "False: In Java, if you name a variable with all uppercase letters and initialize it to some value, attempting to subsequently change its value would not lead to a compile-time error, but it is not a good naming convention.", "Correct", 1)
    ,@(6, "first-year", 0, "False: In Java, after you create and initialize an int[] array, you can change its length by using the method Array.resize().", "Incorrect", 1)
    ,@(7, "first-year", 0, "False: In Java, a[0.0] is not a valid expression that gives the first element of an array of type char[], but a[0] does give the first element.", "Correct", 1)
    ,@(8, "first-year", 0, "False: In Java, it is not possible to declare, create, and initialize a String[] array without using the keyword new.", "Incorrect", 1)
    ,@(9, "first-year", 0, "False: In Java, if a[] and b[] are two different arrays of the same type and length, then the expression (a == b) evaluates to false, even if the corresponding array elements are equal.", "Correct", 1)
    ,@(10, "first-year", 0, "True: In Java, the elements of an array of type int[] are stored contiguously in the computer’s memory (i.e., in consecutive memory locations).", "Correct", 1)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
